# "updated annex and cbo transformation"
# - Rename the info-sheet "TRANSF_ANNEX" header (column E) to "TRANSF_ARG1"
#   and add a new "TRANSF_ARG2" header in column F (the single ANNEX
#   argument became two explicit TRANSF args).
# - The CBO/GDP-SPF row's argument value ("GDP") shifts from column E to
#   the new column F to line up with the new TRANSF_ARG1/TRANSF_ARG2
#   headers.
# - The "info" sheet becomes the active/selected sheet instead of "legend".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("info")

# --- give the new column F the same look as columns C:E before we touch it ---
$ws1.Range("E1").Copy()
$ws1.Range("F1").PasteSpecial(-4122)
$ws1.Range("E8").Copy()
$ws1.Range("F8").PasteSpecial(-4122)
$ws1.Columns("F").ColumnWidth = $ws1.Columns("E").ColumnWidth

# --- header row: TRANSF_ANNEX -> TRANSF_ARG1, plus new TRANSF_ARG2 ---
$ws1.Range("E1").Value = "TRANSF_ARG1"
$ws1.Range("F1").Value = "TRANSF_ARG2"

# --- row 8 (GDP SPF / CBO row): argument moves from column E to column F ---
$ws1.Range("F8").Value = $ws1.Range("E8").Value2
$ws1.Range("E8").Clear()

# --- make "info" the active sheet (was "legend") ---
$ws1.Activate()
$ws1.Range("A1").Select()
